$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 84; existing rows 84-131 shift down to 85-132.
$ws.Rows("84:84").Insert()

$ws.Range("A84").Value = 2
$ws.Range("B84").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44651
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100103
$ws.Range("H84").Value = "Frutos de hueso (carozo)"
$ws.Range("I84").Value = 100103004
$ws.Range("J84").Value = "Durazno"
$ws.Range("K84").Value = "Kakamas"
$ws.Range("L84").Value = "Primera"
$ws.Range("M84").Value = 16
$ws.Range("N84").Value = 450000
$ws.Range("O84").Value = 460000
$ws.Range("P84").Value = 455000
$ws.Range("Q84").Value = "$/bins (400 kilos)"
$ws.Range("R84").Value = "Región de O'Higgins"
$ws.Range("S84").Value = 1138
$ws.Range("T84").Value = 400
